# CT - Manter Usuário.xlsx
# Inclusão de novos casos de teste
# Incluído: Manter Sprint;
# Alterados: Manter Empresa e Manter Usuário.

$wb = $excel.ActiveWorkbook

$wsCapa     = $wb.Worksheets.Item(1)   # Capa
$wsVersao   = $wb.Worksheets.Item(2)   # Versão
$wsScripts  = $wb.Worksheets.Item(3)   # Scripts

# --- Sheet "Capa": give the subtitle row (row 5) a bit more height ---
$wsCapa.Rows.Item(5).RowHeight = 18

# --- Sheet "Versão": add a new version-history entry (row 5) ---
$wsVersao.Range("A5").Value = "04.00 - 18/04/2014"
$wsVersao.Range("B5").Value = "Rodrigo Melo"
$wsVersao.Rows.Item(5).RowHeight = 15.75

# --- Sheet "Scripts": add new test case #13 (row 14) ---
$wsScripts.Range("A14").Value = 13
$wsScripts.Range("B14").Value = "Efetuar a inclusão de um mesmo Usuário em mais de uma empresa."
